$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-MarkCellText {
    param($table, $rowIndex, $colIndex, $newText)
    $cell = $table.Cell($rowIndex, $colIndex)
    $r = $cell.Range
    # The cell range's last "character" is the end-of-cell marker, so only
    # touch the single visible character that holds the mark - this keeps
    # the run's existing character formatting (sz/szCs/lang/etc.) intact.
    $sub = $d.Range($r.Start, $r.Start + 1)
    $sub.Text = $newText
}

# Grading-mark column is column 4. Row numbers below are the 1-based table
# row index (row 1 is the header row).
# Row 2  (label "1"):  x -> X
Set-MarkCellText $t 2 4 "X"
# Row 3  (label "2"):  x -> X
Set-MarkCellText $t 3 4 "X"
# Row 4  (label "3"):  x -> X
Set-MarkCellText $t 4 4 "X"
# Row 6  (label "13"): x -> X
Set-MarkCellText $t 6 4 "X"
# Row 11 (label "28"): x -> X
Set-MarkCellText $t 11 4 "X"
# Row 12 (label "32"): ? -> X
Set-MarkCellText $t 12 4 "X"
# Row 14 (label "36"): ? -> X
Set-MarkCellText $t 14 4 "X"

# Row 9 (label "24"): the mark stays "?" but its run is bumped up a notch -
# sz 28 -> 32 half-points (14pt -> 16pt) and szCs 44 -> 48 (22pt -> 24pt).
# Font.Size only ever touches <w:sz>, never <w:szCs>, in this COM host, so
# rewrite the paragraph's own OOXML (via InsertXML) to land both values -
# every other attribute/element is carried over unchanged.
$cell9 = $t.Cell(9, 4)
$para9 = $cell9.Range.Paragraphs.Item(1)
$pr9 = $para9.Range

$paraXml = '<w:p w14:paraId="43EA636D" w14:textId="0C19CE14" w:rsidR="00D9378E" w:rsidRPr="000B3668" w:rsidRDefault="00A2747C" w:rsidP="0036689A">' +
           '<w:pPr><w:jc w:val="both"/><w:rPr><w:noProof/><w:sz w:val="12"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
           '<w:r w:rsidRPr="00A2747C"><w:rPr><w:noProof/><w:sz w:val="32"/><w:szCs w:val="48"/><w:lang w:val="en-GB"/></w:rPr><w:t>?</w:t></w:r>' +
           '</w:p>'

$xmlFrag = '<?xml version="1.0"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' + $paraXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

$pr9.InsertXML($xmlFrag) | Out-Null
